# Fixed #348 Hyperlinks from sub-template does not work.
#
# The "4", "2" and "finished" paragraphs each start with a stray empty
# run (<w:r><w:rPr>...</w:rPr><w:t/></w:r>) left over from a previous
# edit, and the "finished" paragraph also carries a stray trailing
# empty run (<w:r><w:t/></w:r>). These are removed, and the "finished"
# paragraph gets a <w:proofErr w:type="gramStart"/><w:proofErr
# w:type="gramEnd"/> pair inserted before its run, matching the other
# paragraphs in the document that already bracket their first sentence
# with proofErr markers.
#
# We rebuild each of these three paragraphs from exact OOXML so the
# stray empty runs are dropped while every other attribute (rsids,
# language run properties, etc.) is preserved unchanged.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$paraWith4 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00735354" w:rsidRPr="00DC5685" w:rsidRDefault="00735354" w:rsidP="00F5495F"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>4</w:t></w:r></w:p>'

$paraWith2 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00795F22" w:rsidRDefault="00795F22" w:rsidP="00735354"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r></w:p>'

$paraFinished = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00795F22" w:rsidRDefault="00795F22" w:rsidP="00735354"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>finished</w:t></w:r></w:p>'

# Paragraphs: 1 = intro sentence, 2 = "4", 3 = "2", 4 = "finished",
# 5 = closing sentence. Replace each target paragraph's OOXML in place
# (only once we have confirmed we are looking at the expected paragraph;
# Paragraph.Range.Text includes the trailing paragraph mark, hence
# "-like" with a trailing wildcard instead of an exact "-eq" match).
$p4 = $d.Paragraphs.Item(2)
if ($p4.Range.Text -like "4*") {
    [void]$p4.Range.InsertXML($paraWith4)
}

$p2 = $d.Paragraphs.Item(3)
if ($p2.Range.Text -like "2*") {
    [void]$p2.Range.InsertXML($paraWith2)
}

$pFinished = $d.Paragraphs.Item(4)
if ($pFinished.Range.Text -like "finished*") {
    [void]$pFinished.Range.InsertXML($paraFinished)
}
